$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.612.77'
$ws.Cells.Item(2, 5).Value = '  +4.48%  '
$ws.Cells.Item(3, 4).Value = '2.266.31'
$ws.Cells.Item(3, 5).Value = '  +1.94%  '
$ws.Cells.Item(4, 4).Value = '''1.00'
$ws.Cells.Item(4, 5).Value = '  +0.34%  '
$ws.Cells.Item(5, 4).Value = '''232.06'
$ws.Cells.Item(5, 5).Value = '  +0.63%  '
$ws.Cells.Item(6, 4).Value = '''0.625'
$ws.Cells.Item(6, 5).Value = '  +0.05%  '
$ws.Cells.Item(7, 4).Value = '''61.18'
$ws.Cells.Item(7, 5).Value = '  +0.29%  '
$ws.Cells.Item(8, 4).Value = '''1.00'
$ws.Cells.Item(8, 5).Value = '  -0.04%  '
$ws.Cells.Item(9, 5).Value = '  +2.82%  '
$ws.Cells.Item(10, 4).Value = '''0.0918'
$ws.Cells.Item(10, 5).Value = '  +3.32%  '
$ws.Cells.Item(11, 5).Value = '  +0.51%  '
$ws.Cells.Item(12, 4).Value = '2.602.16'
$ws.Cells.Item(12, 5).Value = '  +1.89%  '
$ws.Cells.Item(13, 4).Value = '''15.71'
$ws.Cells.Item(13, 5).Value = '  +0.47%  '
$ws.Cells.Item(14, 4).Value = '''22.60'
$ws.Cells.Item(14, 5).Value = '  +4.20%  '
$ws.Cells.Item(15, 5).Value = '  +2.83%  '
$ws.Cells.Item(16, 4).Value = '''0.809'
$ws.Cells.Item(16, 5).Value = '  +1.63%  '
$ws.Cells.Item(17, 4).Value = '2.264.71'
$ws.Cells.Item(17, 5).Value = '  +1.88%  '
$ws.Cells.Item(18, 4).Value = '43.445.10'
$ws.Cells.Item(18, 5).Value = '  +4.48%  '
$ws.Cells.Item(19, 4).Value = '0.0₃0934'
$ws.Cells.Item(19, 5).Value = '  +4.43%  '
$ws.Cells.Item(20, 2).Value = 'Litecoin'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(20, 4).Value = '''73.02'
$ws.Cells.Item(20, 5).Value = '  +0.42%  '
$ws.Cells.Item(21, 2).Value = 'Uniswap'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(21, 4).Value = '''6.21'
$ws.Cells.Item(21, 5).Value = '  +2.86%  '
$ws.Cells.Item(22, 4).Value = '''248.39'
$ws.Cells.Item(22, 5).Value = '  -0.52%  '
$ws.Cells.Item(23, 5).Value = '  +8.54%  '
$ws.Cells.Item(24, 4).Value = '''0.999'
$ws.Cells.Item(24, 5).Value = '  -0.09%  '
$ws.Cells.Item(25, 5).Value = '  +4.94%  '
$ws.Cells.Item(26, 4).Value = '''9.80'
$ws.Cells.Item(26, 5).Value = '  +2.36%  '
$ws.Cells.Item(27, 2).Value = 'Monero'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(27, 4).Value = '''169.71'
$ws.Cells.Item(27, 5).Value = '  +1.28%  '
$ws.Cells.Item(28, 2).Value = 'Kaspa'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(28, 4).Value = '''0.143'
$ws.Cells.Item(28, 5).Value = '  +1.90%  '
$ws.Cells.Item(29, 5).Value = '  +5.94%  '
$ws.Cells.Item(30, 4).Value = '''20.59'
$ws.Cells.Item(30, 5).Value = '  +3.30%  '
$ws.Cells.Item(31, 4).Value = '''2.69'
$ws.Cells.Item(31, 5).Value = '  +2.48%  '
$ws.Cells.Item(32, 5).Value = '  -0.85%  '
$ws.Cells.Item(33, 5).Value = '  +2.77%  '
$ws.Cells.Item(34, 5).Value = '  +3.06%  '
$ws.Cells.Item(35, 4).Value = '''0.0656'
$ws.Cells.Item(35, 5).Value = '  +5.55%  '
$ws.Cells.Item(36, 4).Value = '''6.46'
$ws.Cells.Item(36, 5).Value = '  -1.40%  '
$ws.Cells.Item(37, 4).Value = '''2.40'
$ws.Cells.Item(37, 5).Value = '  +2.23%  '
$ws.Cells.Item(38, 4).Value = '''3.60'
$ws.Cells.Item(38, 5).Value = '  -1.80%  '
$ws.Cells.Item(40, 4).Value = '''1.00'
$ws.Cells.Item(40, 5).Value = '  +0.45%  '
$ws.Cells.Item(41, 2).Value = 'FraxShare'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(41, 4).Value = '''8.70'
$ws.Cells.Item(41, 5).Value = '  +1.57%  '
$ws.Cells.Item(42, 2).Value = 'TerraClassic'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Cells.Item(42, 4).Value = '''0.000222'
$ws.Cells.Item(42, 5).Value = '  -9.76%  '
$ws.Cells.Item(43, 4).Value = '''0.0973'
$ws.Cells.Item(43, 5).Value = '  -0.61%  '
$ws.Cells.Item(44, 4).Value = '''1.21'
$ws.Cells.Item(44, 5).Value = '  +0.23%  '
$ws.Cells.Item(45, 2).Value = 'Aave'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(45, 4).Value = '''97.28'
$ws.Cells.Item(45, 5).Value = '  -1.46%  '
$ws.Cells.Item(46, 2).Value = 'FTXToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Cells.Item(46, 4).Value = '''4.40'
$ws.Cells.Item(46, 5).Value = '  -9.59%  '
$ws.Cells.Item(47, 4).Value = '1.468.42'
$ws.Cells.Item(47, 5).Value = '  +0.30%  '
$ws.Cells.Item(48, 4).Value = '''16.75'
$ws.Cells.Item(48, 5).Value = '  +1.64%  '
$ws.Cells.Item(49, 4).Value = '''1.08'
$ws.Cells.Item(49, 5).Value = '  +0.74%  '
$ws.Cells.Item(50, 5).Value = '  -1.43%  '
$ws.Cells.Item(51, 4).Value = '''2.23'
$ws.Cells.Item(51, 5).Value = '  +6.24%  '
